$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections on existing cells (rewording / punctuation fixes) ---
$ws.Range("D42").Value = "Once Above high, same Above high indicator apply"
$ws.Range("D43").Value = "If next candles closes below mid wicks consider and nearest another level consider w wicks"
$ws.Range("D44").Value = "Market resist from high again"
$ws.Range("C82").Value = "Ambush if strong red and no other level nearby,"

# --- Insert a new data row before the final (blank) footer row ---
# Push the old blank footer row (currently row 85) down to row 86,
# copying its formatting exactly so no new style entries are created.
$ws.Range("A85:D85").Copy()
$ws.Range("A86:D86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Make row 85 a new data row, re-using the formatting from row 84
# (the row directly above it, part of the same "Weak" group).
$ws.Range("A84:D84").Copy()
$ws.Range("A85:D85").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Ensure the pushed-down footer row keeps its standard row height.
$ws.Rows.Item(86).RowHeight = 15

# New row's content.
$ws.Range("C85").Value = "If 1st candle Touched yest mid and close below"

# --- Update the AutoFilter range to cover the new last row (86) ---
$ws.Range("A1:D86").AutoFilter()
$ws.Range("A1:D86").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync ---
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$D`$86"

# --- Restore the view's scroll position / selection to the new last rows ---
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C86").Select()
